$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: shift existing rows down to their new positions.
# Using Range.Copy(destination) (restricted to columns A:D) so that the
# exact existing cell style indices are reused instead of new ones being
# fabricated. Must proceed bottom-up so a source row is always read
# before anything is written into it.

$ws.Range("A17:D17").Copy($ws.Range("A20:D20"))
$ws.Range("A16:D16").Copy($ws.Range("A18:D18"))
$ws.Range("A15:D15").Copy($ws.Range("A19:D19"))
$ws.Range("A15:D15").Copy($ws.Range("A17:D17"))
$ws.Range("A14:D14").Copy($ws.Range("A16:D16"))
$ws.Range("A13:D13").Copy($ws.Range("A15:D15"))
$ws.Range("A12:D12").Copy($ws.Range("A14:D14"))
$ws.Range("A11:D11").Copy($ws.Range("A13:D13"))
$ws.Range("A10:D10").Copy($ws.Range("A12:D12"))
$ws.Range("A9:D9").Copy($ws.Range("A11:D11"))

$excel.CutCopyMode = 0

# --- Step 2: write the new/changed cell values ---
# The brand-new strings are written to column D first (in the order the
# original workbook first introduces them in its shared-string table)
# so the resulting shared-string table ordering matches.

$ws.Cells.Item(15,4).Value = "glove says don't show turn signal"
$ws.Cells.Item(16,4).Value = "glove says don't show turn signal"
$ws.Cells.Item(18,4).Value = "phone says don't show stop sign"
$ws.Cells.Item(9,4).Value  = "nav says don't blink leds"
$ws.Cells.Item(10,4).Value = "nav says don't blink leds"
$ws.Cells.Item(11,4).Value = "disconnect"
$ws.Cells.Item(12,4).Value = "disconnect"
$ws.Cells.Item(19,4).Value = "disconnect"

# Row 9/10: L/R P 4 "nav says don't blink leds"
$ws.Cells.Item(9,1).Value = "L"
$ws.Cells.Item(9,2).Value = "P"
$ws.Cells.Item(9,3).Value = 4

$ws.Cells.Item(10,1).Value = "R"
$ws.Cells.Item(10,2).Value = "P"
$ws.Cells.Item(10,3).Value = 4

# Row 11/12: L/R P 5 "disconnect"
$ws.Cells.Item(11,1).Value = "L"
$ws.Cells.Item(11,2).Value = "P"
$ws.Cells.Item(11,3).Value = 5

$ws.Cells.Item(12,1).Value = "R"
$ws.Cells.Item(12,2).Value = "P"
$ws.Cells.Item(12,3).Value = 5

# Row 13/14: B L 1 glove says left / B R 2 glove says right -- unchanged, already shifted correctly

# Row 15/16: B L/R 4 "glove says don't show turn signal"
$ws.Cells.Item(15,1).Value = "B"
$ws.Cells.Item(15,2).Value = "L"
$ws.Cells.Item(15,3).Value = 4

$ws.Cells.Item(16,1).Value = "B"
$ws.Cells.Item(16,2).Value = "R"
$ws.Cells.Item(16,3).Value = 4

# Row 17: B P 3 "phone says stop"  (count changes from 1 to 3)
$ws.Cells.Item(17,1).Value = "B"
$ws.Cells.Item(17,2).Value = "P"
$ws.Cells.Item(17,3).Value = 3
$ws.Cells.Item(17,4).Value = "phone says stop"

# Row 18: B P 4 "phone says don't show stop sign"
$ws.Cells.Item(18,1).Value = "B"
$ws.Cells.Item(18,2).Value = "P"
$ws.Cells.Item(18,3).Value = 4

# Row 19: B P 5 "disconnect"
$ws.Cells.Item(19,1).Value = "B"
$ws.Cells.Item(19,2).Value = "P"
$ws.Cells.Item(19,3).Value = 5

# Row 20: (any)(any)9 ERROR -- already shifted correctly, values unchanged

$ws.Range("D19").Select()
